# Add files via upload
# - Fixes a typo: "Rdrigo" -> "Rodrigo" for the user in row 27
# - Appends a new data row (row 28) for user "Teresa"
# - Updates the active selection to match the saved view state (F30)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix misspelled name in existing row
$ws.Range("A27").Value = "Rodrigo"

# New row of data appended to the table
$ws.Range("A28").Value = "Teresa"
$ws.Range("B28").Value = 9012
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 3240
$ws.Range("E28").Value = "Saramago"
$ws.Range("F28").Value = "Alameda"
$ws.Range("G28").Value = "38.7401, -9.1340"
$ws.Range("H28").Value = 46
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = 28

# Match the saved selection state from the workbook
$ws.Range("F30").Select()
